# Rename ObjTables document/table metadata attributes to lowerCamelCase.
# Affects the inline "!!ObjTables ..." / "!!!ObjTables ..." header strings
# stored in cell A1 (and A2 on the table-of-contents sheet) of each sheet.

$wb = $excel.ActiveWorkbook

$tocSheet    = $wb.Worksheets.Item("!!_Table of contents")
$schemaSheet = $wb.Worksheets.Item("!!_Schema")
$parentSheet = $wb.Worksheets.Item("!!Parent")
$childSheet  = $wb.Worksheets.Item("!!Child")

$tocSheet.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$tocSheet.Range("A2").Value = "!!ObjTables type='TableOfContents' description='Table/model and column/attribute definitions' date='2019-09-18 13:17:59' objTablesVersion='2.0'"

$schemaSheet.Range("A1").Value = "!!ObjTables type='Schema' description='Table/model and column/attribute definitions' objTablesVersion='0.0.8'"

$parentSheet.Range("A1").Value = "!!ObjTables type='Data' id='Parent' name='Parent' date='2019-09-18 13:17:59' objTablesVersion='0.0.8'"

$childSheet.Range("A1").Value = "!!ObjTables type='Data' id='Child' name='Child' date='2019-09-18 13:17:59' objTablesVersion='0.0.8'"
